# Increment the "Förändrad" (Changed) date column (C) by one day
# for every data row (rows 2 through 267) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 267
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    $cell.Value2 = $cell.Value2 + 1
}
